# Update Name of Algo
# Apply updated KNN-imputed values to the result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "E3"   = 16.664
    "D7"   = -7.535000000000001
    "C8"   = -12.748
    "C10"  = -12.937
    "C12"  = -11.207
    "D15"  = -8.378000000000002
    "C18"  = -13.806
    "D18"  = -8.352
    "E18"  = 16.367
    "E19"  = 16.606
    "D20"  = -7.517
    "E27"  = 16.512
    "D29"  = -7.292
    "D30"  = -7.272
    "D31"  = -7.606
    "E31"  = 17.092
    "C37"  = -13.475
    "E38"  = 16.726
    "D40"  = -7.641
    "E42"  = 16.555
    "E44"  = 16.825
    "E47"  = 16.471
    "D50"  = -8.008000000000001
    "C55"  = -13.958
    "E58"  = 16.701
    "E65"  = 17.119
    "C68"  = -11.232
    "D68"  = -7.053
    "E73"  = 16.667
    "D76"  = -7.312
    "C77"  = -13.117
    "C78"  = -13.214
    "C81"  = -13.234
    "C82"  = -11.737
    "D87"  = -8.339
    "D88"  = -8.177000000000001
    "E90"  = 16.492
    "E94"  = 17.741
    "E95"  = 17.397
    "D96"  = -7.267
    "D98"  = -8.242000000000001
    "D101" = -7.616
    "E101" = 16.44
    "D102" = -8.036
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$wb.Save()
